$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 557, shifting rows 557:653 down to 558:654
$ws.Rows("557:557").Insert()

# Populate the new row 557 with the new data record
$ws.Range("A557").Value = 5
$ws.Range("B557").Value = "Macroferia Regional de Talca"
$ws.Range("C557").Value = "Maule"
$ws.Range("D557").Value = 44522
$ws.Range("D557").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E557").Value = 7
$ws.Range("F557").Value = "Fruta"
$ws.Range("G557").Value = 100103
$ws.Range("H557").Value = "Frutos de hueso (carozo)"
$ws.Range("I557").Value = 100103004
$ws.Range("J557").Value = "Durazno"
$ws.Range("K557").Value = "Florida King"
$ws.Range("L557").Value = "Primera"
$ws.Range("M557").Value = 50
$ws.Range("N557").Value = 24000
$ws.Range("O557").Value = 24000
$ws.Range("P557").Value = 24000
$ws.Range("Q557").Value = "$/bandeja 15 kilos granel"
$ws.Range("R557").Value = "Provincia de San Felipe de Aconcagua"
$ws.Range("S557").Value = 1600
$ws.Range("T557").Value = 15
